$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 48: "Stock adjustment journals" entry updated with new progress figures
$ws.Range("D48").NumberFormat = $ws.Range("D47").NumberFormat
$ws.Range("D48").Value = 0.95
$ws.Range("E48").Value = 1
$ws.Range("H48").Value = 1
$ws.Range("J48").Value = "Tiem reduced from 8 hours"

# Leave the selection / scroll state where the author left it when saving
[void]$ws.Range("D49").Select()
